{"js": "// The document ends with a signature block listing attendee names, each\n// in its own bold paragraph, followed by one extra trailing empty\n// paragraph (also formatted bold). That trailing empty paragraph is\n// removed so the document now ends right after the last name\n// (\"Swaroop Dattatraya Patil:\").\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The very last paragraph in the body is the empty one we need to drop.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.delete();\n\nawait context.sync();\n", "ps1": "# The document's signature block ends with one paragraph per attendee\n# name (e.g. \"Swaroop Dattatraya Patil:\") followed by one extra, empty\n# trailing paragraph. Remove that trailing empty paragraph so the\n# document now ends right after the last attendee's name.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($count)\n$previousPara = $d.Paragraphs.Item($count - 1)\n\n# Build a range that starts right at the end of the previous paragraph's\n# own text (i.e. just before its paragraph mark) and runs through the\n# end of the very last (empty) paragraph. Deleting it removes only the\n# trailing empty paragraph's mark while leaving the previous paragraph,\n# its own mark, and its formatting completely untouched.\n$rangeToRemove = $d.Range($previousPara.Range.End - 1, $lastPara.Range.End)\n$rangeToRemove.Delete()\n"}
